$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-03-13 Wednesday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-03-14 Thursday", 2)
$d.Content.Find.Execute("971×2=1942", $true, $false, $false, $false, $false, $true, 1, $false, "272×7=1904", 2)
$d.Content.Find.Execute("569×9=5121", $true, $false, $false, $false, $false, $true, 1, $false, "405×3=1215", 2)
$d.Content.Find.Execute("732×4=2928", $true, $false, $false, $false, $false, $true, 1, $false, "534×9=4806", 2)
$d.Content.Find.Execute("851×9=7659", $true, $false, $false, $false, $false, $true, 1, $false, "172×6=1032", 2)
$d.Content.Find.Execute("830×8=6640", $true, $false, $false, $false, $false, $true, 1, $false, "501×3=1503", 2)
$d.Content.Find.Execute("936×8=7488", $true, $false, $false, $false, $false, $true, 1, $false, "375×5=1875", 2)
$d.Content.Find.Execute("345×6=2070", $true, $false, $false, $false, $false, $true, 1, $false, "173×8=1384", 2)
$d.Content.Find.Execute("195×8=1560", $true, $false, $false, $false, $false, $true, 1, $false, "622×9=5598", 2)
$d.Content.Find.Execute("569×6=3414", $true, $false, $false, $false, $false, $true, 1, $false, "909×5=4545", 2)
$d.Content.Find.Execute("390×9=3510", $true, $false, $false, $false, $false, $true, 1, $false, "785×4=3140", 2)
$d.Content.Find.Execute("918×2=1836", $true, $false, $false, $false, $false, $true, 1, $false, "876×2=1752", 2)
$d.Content.Find.Execute("686×6=4116", $true, $false, $false, $false, $false, $true, 1, $false, "164×4=656", 2)
$d.Content.Find.Execute("383×4=1532", $true, $false, $false, $false, $false, $true, 1, $false, "891×6=5346", 2)
$d.Content.Find.Execute("562×6=3372", $true, $false, $false, $false, $false, $true, 1, $false, "452×4=1808", 2)
$d.Content.Find.Execute("247×3=741", $true, $false, $false, $false, $false, $true, 1, $false, "838×6=5028", 2)
$d.Content.Find.Execute("334×3=1002", $true, $false, $false, $false, $false, $true, 1, $false, "975×7=6825", 2)
$d.Content.Find.Execute("459×8=3672", $true, $false, $false, $false, $false, $true, 1, $false, "744×8=5952", 2)
$d.Content.Find.Execute("898×4=3592", $true, $false, $false, $false, $false, $true, 1, $false, "159×4=636", 2)
$d.Content.Find.Execute("773×3=2319", $true, $false, $false, $false, $false, $true, 1, $false, "794×8=6352", 2)
$d.Content.Find.Execute("246×6=1476", $true, $false, $false, $false, $false, $true, 1, $false, "234×2=468", 2)
$d.Content.Find.Execute("477×8=3816", $true, $false, $false, $false, $false, $true, 1, $false, "941×7=6587", 2)
$d.Content.Find.Execute("830×4=3320", $true, $false, $false, $false, $false, $true, 1, $false, "504×2=1008", 2)
$d.Content.Find.Execute("909×6=5454", $true, $false, $false, $false, $false, $true, 1, $false, "715×9=6435", 2)
$d.Content.Find.Execute("267×4=1068", $true, $false, $false, $false, $false, $true, 1, $false, "595×6=3570", 2)
$d.Content.Find.Execute("730×2=1460", $true, $false, $false, $false, $false, $true, 1, $false, "944×8=7552", 2)

Write-Output "replacements complete"
